$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on Price cells that would otherwise be auto-converted to numbers
$textCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D14", "D15", "D16", "D18", "D19", "D20", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values from the crypto data refresh
$ws.Range("D2").Value = "30.231.90"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").Value = "1.861.65"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "235.81"
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("D6").Value = "0.9992"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").Value = "0.4695"
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("D8").Value = "0.2896"
$ws.Range("D9").Value = "0.06565"
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("D10").Value = "21.84"
$ws.Range("E10").Value = "  +1.93%  "
$ws.Range("D11").Value = "0.07994"
$ws.Range("E11").Value = "  +1.36%  "
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("D13").Value = "1.862.44"
$ws.Range("E13").Value = "  -0.85%  "
$ws.Range("D14").Value = "5.110"
$ws.Range("E14").Value = "  -0.03%  "
$ws.Range("D15").Value = "0.6785"
$ws.Range("E15").Value = "  +0.26%  "
$ws.Range("D16").Value = "268.62"
$ws.Range("E16").Value = "  -3.76%  "
$ws.Range("D17").Value = "30.210.14"
$ws.Range("E17").Value = "  -0.45%  "
$ws.Range("D18").Value = "13.61"
$ws.Range("E18").Value = "  +6.81%  "
$ws.Range("D19").Value = "0.000007647"
$ws.Range("E19").Value = "  +4.60%  "
$ws.Range("D20").Value = "0.9992"
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("D21").Value = "2.104.62"
$ws.Range("E21").Value = "  -0.81%  "
$ws.Range("D22").Value = "0.9998"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").Value = "5.225"
$ws.Range("E23").Value = "  -4.62%  "
$ws.Range("D24").Value = "6.187"
$ws.Range("E24").Value = "  +0.48%  "
$ws.Range("D25").Value = "167.16"
$ws.Range("E25").Value = "  +1.08%  "
$ws.Range("D26").Value = "9.188"
$ws.Range("E26").Value = "  +0.33%  "
$ws.Range("D27").Value = "18.96"
$ws.Range("E27").Value = "  -1.13%  "
$ws.Range("D28").Value = "1.949"
$ws.Range("E28").Value = "  +0.78%  "
$ws.Range("D29").Value = "1.375"
$ws.Range("E29").Value = "  -0.54%  "
$ws.Range("D30").Value = "0.09906"
$ws.Range("E30").Value = "  +2.22%  "
$ws.Range("D31").Value = "4.338"
$ws.Range("E31").Value = "  -1.35%  "
$ws.Range("D32").Value = "1.464"
$ws.Range("E32").Value = "  -0.75%  "
$ws.Range("D33").Value = "4.042"
$ws.Range("E33").Value = "  -1.54%  "
$ws.Range("D34").Value = "0.04712"
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("D35").Value = "1.123"
$ws.Range("E35").Value = "  -0.40%  "
$ws.Range("D36").Value = "0.7017"
$ws.Range("E36").Value = "  -0.69%  "
$ws.Range("D37").Value = "2.703"
$ws.Range("E37").Value = "  -0.84%  "
$ws.Range("D38").Value = "0.01871"
$ws.Range("E38").Value = "  +0.53%  "
$ws.Range("D39").Value = "2.609"
$ws.Range("E39").Value = "  +2.78%  "
$ws.Range("D40").Value = "6.331"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").Value = "73.49"
$ws.Range("E41").Value = "  -1.31%  "
$ws.Range("D42").Value = "1.938"
$ws.Range("E42").Value = "  -0.91%  "
$ws.Range("D43").Value = "0.8393"
$ws.Range("E43").Value = "  -1.35%  "
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").Value = "0.9984"
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "103.67"
$ws.Range("E45").Value = "  -0.16%  "
$ws.Range("B46").Value = "TheSandbox"
$ws.Range("C46").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D46").Value = "0.4148"
$ws.Range("E46").Value = "  -0.99%  "
$ws.Range("D47").Value = "9.192"
$ws.Range("E47").Value = "  -0.78%  "
$ws.Range("D48").Value = "7.064"
$ws.Range("E48").Value = "  -2.18%  "
$ws.Range("D49").Value = "935.23"
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("D50").Value = "34.08"
$ws.Range("E50").Value = "  -0.47%  "
$ws.Range("D51").Value = "0.05659"
$ws.Range("E51").Value = "  +0.44%  "
